# Update EUR->ARS rate: add a new data row (2025-10-01T21:21:00Z) to the
# quote history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 52

# Write the three new values as text-literal formulas first (wrapping in
# quotes forces Excel to treat them as strings rather than auto-converting
# "2025-10-01" / "21:21:00" into date/time serial numbers), then convert
# the formulas in place to plain static values so the cells end up as
# ordinary text cells (matching the existing rows), with no formulas left
# behind and no extra number-format/style entries introduced.
$ws.Cells.Item($newRow, 1).Formula = "=""2025-10-01"""
$ws.Cells.Item($newRow, 2).Formula = "=""21:21:00"""
$ws.Cells.Item($newRow, 3).Formula = "=""1.00 EUR = 1,771.6830"""

$rowRange = $ws.Range("A" + $newRow + ":C" + $newRow)
$rowRange.Copy()
$rowRange.PasteSpecial(-4163)  # xlPasteValues
